# Applies the "added temple and forest enemy concepts" edit:
#  1. Colors everything in the "-All main areas and the goals in that
#     Area(finding weapon, defeating boss)" bullet green (RGB 00B050),
#     except the leading "-".
#  2. Normalizes a couple of runs that had been split mid-word back into
#     single runs ("Mecha"+"nics(" -> "Mechanics(", and
#     "-Various assets fo"+"r each Area" -> "-Various assets for each Area")
#     by doing an identity Find/Replace, which causes same-formatted
#     adjacent runs to coalesce.

$d = $word.ActiveDocument

# --- 1. Colorize the "-All main areas..." bullet (skip the leading "-") ---
$r = $d.Content
$found = $r.Find.Execute(
    "All main areas and the goals in that Area(finding weapon, defeating boss)",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    # wdColor values are packed 0xBBGGRR, so RGB 00B050 -> 0x0050B000
    $r.Font.Color = 0x0050B000
}

# --- 2. Re-merge "Mecha" / "nics(" into a single "Mechanics(" run ---
$d.Content.Find.Execute("Mecha", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Mecha", 2)

# --- 3. Re-merge "-Various assets fo" / "r each Area" into one run ---
$d.Content.Find.Execute("-Various assets fo", $true, $false, $false, $false,
                         $false, $true, 1, $false, "-Various assets fo", 2)
